# Adds the missing body text to the two placeholder-only slides
# (sldId 287 -> slide index 3, "Project design"; sldId 290 -> slide
# index 7, "Response variables") per the commit "tilfojet to slides
# til PPT".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 3 ("Project design") - 3 empty body placeholders get filled.
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Shape id=3 "Pladsholder til tekst 2" -> agenda item 1
$tr = $s3.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "1)" + "`r" + "Basics of Machine Learning" + "`r" + "What is Deep learning?"

# Shape id=4 "Pladsholder til tekst 3" -> agenda item 2
$tr = $s3.Shapes.Item(3).TextFrame.TextRange
$tr.Text = "2)" + "`r" + "Data sources" + "`r" + "Data collection" + "`r" + "Data manipulation"

# Shape id=5 "Pladsholder til tekst 4" -> agenda item 3
$tr = $s3.Shapes.Item(4).TextFrame.TextRange
$tr.Text = "3)" + "`r" + "Deep learning model" + "`r" + "Benchmarking" + "`r" + "Results"

# ---------------------------------------------------------------
# Slide 7 ("Response variables") - 2 empty body placeholders get filled.
# ---------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# Shape id=3 "Pladsholder til tekst 2" -> bold intro line + explanatory bullets
$tr = $s7.Shapes.Item(2).TextFrame.TextRange
$introLine = "To forecast electricity prices in period t, we use:"
$tr.Text = $introLine + "`r" + "Prices: t-3, t-2, t-1" + "`r" + "Sales: t-3, t-2, t-1" + "`r" + "Weather: t-1" + "`r" + "GDP: t-1"
$tr.Characters(1, $introLine.Length).Font.Bold = $true

# Shape id=4 "Pladsholder til tekst 3" -> descriptive-statistics bullets
$tr = $s7.Shapes.Item(3).TextFrame.TextRange
$tr.Text = "Sales variance, skewness and curtosis" + "`r" + "Price variance, skewness and curtosis" + "`r" + "Population in 2013" + "`r" + "Longtitude/latitude" + "`r"
